# Add new budget category "Malysh" (Baby) as a new column pair (data + gutter)
# right before the existing "Dom, kvartira" category, i.e. at column P.
# Everything from P onward shifts two columns to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new columns at P:Q (shifts old P.. onward to R..)
$ws.Columns("P:Q").Insert()

# 2) New header cell P1 - copy formatting from the neighbouring header (N1)
#    so it gets the same bold/border/centered style, then set its text.
$ws.Range("N1").Copy()
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("P1").Value = "Малыш"

# 3) Column widths that differ from a plain shift:
#    (ColumnWidth input is offset by -5/6 vs. the stored character width
#    because of Excel's pixel-snapping of the "characters" unit.)
$ws.Columns("F:F").ColumnWidth = 9.166666666666666    # Другое+ : 20 -> 10
$ws.Columns("N:N").ColumnWidth = 19.166666666666668   # Здоровье : 10 -> 20
$ws.Columns("P:P").ColumnWidth = 7.166666666666667    # Малыш (new) : -> 8
$ws.Columns("Q:Q").ColumnWidth = 5.166666666666667    # gutter (new) : -> 6
$ws.Columns("AB:AB").ColumnWidth = 8.166666666666666  # Хоз. Товары : 10 -> 9

# 4) Updated data values (row 2 = Доходы, row 3 = Расходы, row 4 = Итого)
$ws.Range("B2").Value = 137725.5
$ws.Range("D2").Value = 46109.09
$ws.Range("F2").Value = 21263.36
$ws.Range("P2").Value = 0
$ws.Range("AP2").Value = 0

$ws.Range("B3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("J3").Value = 820
$ws.Range("L3").Value = 5678
$ws.Range("N3").Value = 1745.18
$ws.Range("P3").Value = 3785.1
$ws.Range("R3").Value = 5302.15
$ws.Range("T3").Value = 460
$ws.Range("V3").Value = 457
$ws.Range("X3").Value = 6416.26
$ws.Range("Z3").Value = 25023.52
$ws.Range("AB3").Value = 14339
$ws.Range("AD3").Value = 10487
$ws.Range("AF3").Value = 1138
$ws.Range("AH3").Value = 10039.91
$ws.Range("AJ3").Value = 5159
$ws.Range("AL3").Value = 3736.05
$ws.Range("AN3").Value = 0
$ws.Range("AP3").Value = 0

$ws.Range("B4").Value = 137725.5
$ws.Range("D4").Value = 46109.09
$ws.Range("F4").Value = 21263.36
$ws.Range("J4").Value = 820
$ws.Range("L4").Value = 5678
$ws.Range("N4").Value = 1745.18
$ws.Range("P4").Value = 3785.1
$ws.Range("R4").Value = 5302.15
$ws.Range("T4").Value = 460
$ws.Range("V4").Value = 457
$ws.Range("X4").Value = 6416.26
$ws.Range("Z4").Value = 25023.52
$ws.Range("AB4").Value = 14339
$ws.Range("AD4").Value = 10487
$ws.Range("AF4").Value = 1138
$ws.Range("AH4").Value = 10039.91
$ws.Range("AJ4").Value = 5159
$ws.Range("AL4").Value = 3736.05
$ws.Range("AN4").Value = 0
$ws.Range("AP4").Value = 0
